# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" worksheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# Row -> new F value, for the "展览" worksheet
$exhibitionUpdates = @{
    2  = 276
    6  = 1084
    7  = 1424
    9  = 105
    12 = 150
    15 = 1338
    20 = 646
    21 = 33
    22 = 210
    23 = 19
    24 = 5754
    29 = 14342
    30 = 1427
    31 = 200
    32 = 98
    34 = 4437
    35 = 602
    36 = 4188
    37 = 130
}

# Row -> new F value, for the "全部类型" worksheet
$allTypesUpdates = @{
    2  = 276
    6  = 1084
    7  = 1424
    9  = 105
    12 = 150
    15 = 1338
    21 = 646
    23 = 33
    24 = 210
    25 = 19
    27 = 5754
    32 = 14342
    33 = 1427
    34 = 200
    35 = 98
    37 = 4437
    38 = 602
    39 = 4188
    40 = 130
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
